# Book1.xlsx - "Updated unstable to failure"
#
# The scenario list in column A gets:
#   - "E2E_Intra" (row 3) removed entirely (rows below shift up)
#   - the old "E2E_StockRotationReturnDelivery" entry renamed to
#     "E2E_BTS_Bundle" and re-styled to match the earlier "unstable" rows
#   - two new scenario rows appended at the bottom:
#       "E2E_StockRotationReturn" and "E2E_LocalCurr"
#   - the sheet's selected cell moves to D9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# 1) Remove the "E2E_Intra" row (row 3); everything below shifts up one row.
$ws.Rows(3).Delete()

# 2) The row that used to hold "E2E_StockRotationReturnDelivery" is now A9.
#    Rename it to "E2E_BTS_Bundle" and match the formatting used higher up
#    in the list (same look as A2/A3, style index 2).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A9").Value = "E2E_BTS_Bundle"

# 3) Append two new scenario rows at the bottom (A10, A11), matching the
#    formatting of the preceding rows (same look as A8, style index 3).
$ws.Range("A8").Copy() | Out-Null
$ws.Range("A10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A10").Value = "E2E_StockRotationReturn"

$ws.Range("A8").Copy() | Out-Null
$ws.Range("A11").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A11").Value = "E2E_LocalCurr"

$excel.CutCopyMode = $false

# 4) Update the selected cell to D9 to match the saved view state.
$ws.Range("D9").Select() | Out-Null
